$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected Price/Volume columns to Text format before writing so that
# numeric-looking strings (e.g. "0.994") are stored as text, matching the original
# inline-string cell type instead of being coerced into numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.848.42'
$ws.Range('E2').Value = '  +5.01%  '
$ws.Range('D3').Value = '1.610.85'
$ws.Range('E3').Value = '  +3.67%  '
$ws.Range('D5').Value = '213.79'
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('E6').Value = '  +6.88%  '
$ws.Range('D7').Value = '0.994'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('D8').Value = '26.97'
$ws.Range('E8').Value = '  +12.07%  '
$ws.Range('E9').Value = '  +3.19%  '
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('E11').Value = '  +2.54%  '
$ws.Range('D12').Value = '1.842.06'
$ws.Range('E12').Value = '  +3.71%  '
$ws.Range('D13').Value = '1.610.25'
$ws.Range('E13').Value = '  +3.69%  '
$ws.Range('D14').Value = '29.857.72'
$ws.Range('E14').Value = '  +5.00%  '
$ws.Range('D15').Value = '0.538'
$ws.Range('E15').Value = '  +5.59%  '
$ws.Range('E16').Value = '  +3.72%  '
$ws.Range('D17').Value = '244.80'
$ws.Range('E17').Value = '  +6.91%  '
$ws.Range('D18').Value = '63.42'
$ws.Range('E18').Value = '  +3.85%  '
$ws.Range('D19').Value = '7.62'
$ws.Range('E19').Value = '  +3.94%  '
$ws.Range('E20').Value = '  +3.20%  '
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '4.04'
$ws.Range('E22').Value = '  +4.10%  '
$ws.Range('D23').Value = '9.26'
$ws.Range('E23').Value = '  +3.84%  '
$ws.Range('E24').Value = '  +3.90%  '
$ws.Range('D25').Value = '155.85'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('D26').Value = '15.36'
$ws.Range('E26').Value = '  +4.33%  '
$ws.Range('D27').Value = '0.109'
$ws.Range('E27').Value = '  +5.66%  '
$ws.Range('E28').Value = '  +2.73%  '
$ws.Range('D29').Value = '0.995'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +2.78%  '
$ws.Range('D33').Value = '1.442.84'
$ws.Range('E33').Value = '  +4.24%  '
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').Value = '  +3.79%  '
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  +10.22%  '
$ws.Range('E37').Value = '  +2.42%  '
$ws.Range('E39').Value = '  +3.18%  '
$ws.Range('E40').Value = '  +5.38%  '
$ws.Range('D41').Value = '55.61'
$ws.Range('E41').Value = '  +30.35%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('D46').Value = '66.08'
$ws.Range('E46').Value = '  +6.98%  '
$ws.Range('D47').Value = '5.31'
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('D48').Value = '1.753.41'
$ws.Range('E48').Value = '  +3.99%  '
$ws.Range('D49').Value = '87.07'
$ws.Range('E49').Value = '  +2.42%  '
$ws.Range('D50').Value = '0.837'
$ws.Range('E50').Value = '  -4.31%  '
$ws.Range('D51').Value = '0.0₆0102'
$ws.Range('E51').Value = '  +1.46%  '

# Restore the default (Normal) style so no stray text-format styling is left behind.
$dataRange.Style = "Normal"

